$d = $word.ActiveDocument
$p = $d.Paragraphs.Last
$pStart = $p.Range.Start
$pEnd = $p.Range.End
$full = $d.Range($pStart, $pEnd)
$full.Delete()

$endPos = $d.Content.End
$insertRange = $d.Range($endPos, $endPos)
$xml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:ve="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml"><w:body><w:p w:rsidR="00001046" w:rsidRPr="00371DBA" w:rsidRDefault="00371DBA" w:rsidP="00C770D4"><w:pPr><w:pStyle w:val="PlainText"/><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/></w:rPr><w:t>OS X</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:lang w:val="bg-BG"/></w:rPr><w:t xml:space="preserve">. Юнити е основният инструмент за софтуерна обработка на конзолата на Нинтендо </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/></w:rPr><w:t>Wii</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/></w:rPr><w:t xml:space="preserve"> U</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:lang w:val="bg-BG"/></w:rPr><w:t xml:space="preserve"> и има включено безплатно копие към всеки </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/></w:rPr><w:t>Wii</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/></w:rPr><w:t xml:space="preserve"> U </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:lang w:val="bg-BG"/></w:rPr><w:t>лиценз за програмиране.</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:lang w:val="bg-BG"/></w:rPr><w:t xml:space="preserve"> Юнити Технолоджис нарича това добавяне на страничен софтуер „първо в индустрията”. </w:t></w:r></w:p>

    <w:p>
      <w:pPr>
        <w:pStyle w:val="PlainText"/>
        <w:rPr>
          <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:u w:val="single"/>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:u w:val="single"/>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t>Юнити 2</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:lastRenderedPageBreak/>
        <w:t xml:space="preserve">Юнити 2.0 е пуснато в продажба на 11 Октомври 2007г.  по време на първата годишна Юнити конференция </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">(Unite conference). </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t xml:space="preserve">Ключови добавки към съществуващата версия включват картографиращ енджин, мрежова система </w:t>
      </w:r>
      <w:r>
        <w:t>(</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t xml:space="preserve">базирана на </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>RakNet</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">), </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t>динамично добавяне на сенки в реално време, и система за изграждане на игрален юзър интерфейс. В тази версия също е представен Юнити Асет Сървър, добавка която позволява на екипи от програмисти да споделят по лесно придобивките на проекта.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t>На 4 Октомври 2008г. е обявена добавката за публикуване за Айфон. Това позволило на програмистите да напишат игри на Мак и да ги публикуват за Айфон.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t xml:space="preserve">В Юнити 2.5, пуснато на 19 Март 2009г.  е добавен съпорт за писане на игри на Уиндоус. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t xml:space="preserve">През Октомври 2009г. на Юнити конференцията е потвърдено, че Юнити Технолоджис няма да искат вече пари за „инди” версията на Юнити, а вместо това ще я пуснат за свободно ползване. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:u w:val="single"/>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:u w:val="single"/>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t>Юнити 3</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t xml:space="preserve">Юнити 3.0 е пуснато на 4 Октомври 2010г. Новостите включват съпорт за светлинно картографиране с използване на </w:t>
      </w:r>
      <w:r>
        <w:t>Beast</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t xml:space="preserve"> технологията на </w:t>
      </w:r>
      <w:r>
        <w:t>Illuminate Labs</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t xml:space="preserve">, премахване на скрити повърхности предоставено от </w:t>
      </w:r>
      <w:r>
        <w:t>Umbra</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t xml:space="preserve">, обработка на звукови процеси в реално време и поддръжка на </w:t>
      </w:r>
      <w:r>
        <w:t>C# 3.5.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t xml:space="preserve"> Тази версия включва и поддръжка на „визуализация на качеството” за публикуване на игри за Андроид, продавана като отделна добавка. Поддръжката на Андроид официално е обявена на 1 Март 2011г. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t xml:space="preserve">На 10 Ноември 2010г.  е лансиран Юнити Магазин за придобивки, където Юнити потребителите могат онлайн да продават придобивки за своите проекти – произведения на изкуството, кодови системи, аудио и др. – един на друг. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t xml:space="preserve">Юнити 3.4 предостави вградена поддръжка за системата за процесуално структуриране </w:t>
      </w:r>
      <w:r>
        <w:t>“Substance”</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t xml:space="preserve"> на </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Allegorithmic</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t xml:space="preserve">. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t>Юнити 3.5 е реализирано на 14 Февруари 2014 и се отличава с предоставянето на няколко ключови особености доста късно в цикъла на живот на версията. Нова система за частици наречена „Шурикен”,  вградена структура за рутиране и навигиране, детайлно ниво на управление на 3</w:t>
      </w:r>
      <w:r>
        <w:t>D</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t xml:space="preserve"> модели, </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">HDR </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t xml:space="preserve">изпълнение, нови свойства за глобално осветяване, и пренаписване на премахването на скрити повърхности. В тази версия е и представянето на предварителната поддръжка за </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">Adobe Flash </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t xml:space="preserve">и </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">Google Native Client </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t xml:space="preserve">като платформи за публикуване. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:u w:val="single"/>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:u w:val="single"/>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t xml:space="preserve">Юнити 4 </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:lastRenderedPageBreak/>
        <w:t xml:space="preserve">Юнити 4.0 е официално представено на 13 Ноември 2012г.  Главните нови свойства включват нова система за анимиране </w:t>
      </w:r>
      <w:r>
        <w:t>“</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Mecanim</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">”, </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t xml:space="preserve">поддръжка на </w:t>
      </w:r>
      <w:r>
        <w:t>DirectX 11</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t xml:space="preserve">,  и затъмняване в реално време за мобилни платформи. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t>С тази версия компанията обявява смяна на посоката към цикъл на обновление с по-малко нови черти, но за по-кратно време. При това положение следващите версии на 4.Х предоставиха новите особености както следва:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t xml:space="preserve">• Юнити 4.1 представено на 13 Март 2013г. : Профилиране на паметта, поддръжка на </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>AirPlay</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t xml:space="preserve">за </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>iOS</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">, </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t xml:space="preserve">и няколко по-малки ъпдейта за </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Mecanim</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t xml:space="preserve"> и редакция на шейдъри.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t xml:space="preserve">• Юнити 4.2 представено на 22 Юли 2013г. : Поддръжка на </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">Windows Phone 8, Windows Store </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t xml:space="preserve">и </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">BlackBerry </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t>като платформи за публикуване</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">; </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t xml:space="preserve">поддръжка на </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">OpenGL ES 3.0 </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t>за мобилни платформи</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">; </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t xml:space="preserve">вграден съпорт за система за контрол на версиите за </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">Perforce; </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t xml:space="preserve">и възможността да се прекрати процеса на изграждане, когато вече е започнал. </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>(</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t xml:space="preserve">Последното получи аплодисменти при обявяването си на </w:t>
      </w:r>
      <w:r>
        <w:t>Unite Nordic</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t xml:space="preserve"> през 2013г.</w:t>
      </w:r>
      <w:r>
        <w:t>)</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t>• Юнити 4.3 представено 12 Ноември 2013г. : нова 2</w:t>
      </w:r>
      <w:r>
        <w:t>D</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t xml:space="preserve"> структура, включваща поддръжка на 2</w:t>
      </w:r>
      <w:r>
        <w:t>D</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t xml:space="preserve"> представяне и нов 2</w:t>
      </w:r>
      <w:r>
        <w:t>D</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t xml:space="preserve"> физичен енджин </w:t>
      </w:r>
      <w:r>
        <w:t>(</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t xml:space="preserve">предоставен от </w:t>
      </w:r>
      <w:r>
        <w:t>Box2D)</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t xml:space="preserve">• Юнити 4.5 представено на 27 Май 2014г. : не са представени значителни обновления, като вместо това тази версия се фокусира върху оправянето на бъгове, докладвайки за повече от 450 поправки във версията. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t xml:space="preserve">• Юнити 4.6 представено на 26 Ноември 2014г. : нова структура на потребителския интерфейс. Също така версия 4.6.2 представена на 29 Януари 2015г. добавя поддръжка на 64-битови приложения за </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>iOS</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t>На 21 Май 2013г. Главният изпълнителен директор Дейвид Хелгасон обяви че основните версии на добавките за Айфон и Андроид ще бъдат свободно достъпни от Юнити 4.2 нататък.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:u w:val="single"/>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:u w:val="single"/>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t>Юнити 5</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t xml:space="preserve">Юнити 5.0 е пуснато за свободно ползване на 3ти Март 2015г.  като е добавено дългоочакваното глобално осветление в реално време базирано на Геометричната Осветителна технология. Другите основни промени са физически базирани шейдъри, </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">HDR </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t xml:space="preserve">небесни полета, отражателни сонди, нов аудио миксер с ефекти и подобрени условия за анимиране. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t xml:space="preserve">Представена е системата </w:t>
      </w:r>
      <w:r>
        <w:t>Cloud Build (</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t>на цена от 25 долара на месец за физически лица</w:t>
      </w:r>
      <w:r>
        <w:t>)</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t xml:space="preserve"> както и „Доклад за състоянието на играта” и „Анализиране на играта” </w:t>
      </w:r>
      <w:r>
        <w:t>(</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t>също за 25 долара на месец за физицески лица</w:t>
      </w:r>
      <w:r>
        <w:t>)</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t xml:space="preserve">, които записват данните на играчите при реализираните игри – нещо което за много програмисти е било трудна да се направи при Юнити 4.х.  Преди това програмистът е трябвало да пише допълнителен код логвайки се директно в гейм енджина на играча. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:lastRenderedPageBreak/>
        <w:t xml:space="preserve">По-малките подобрения включват: 64-битов едитор за управление на големи проекти, 64-битов съпорт за </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>iOS</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t xml:space="preserve">, ново забавено предаване, графични командни буфери, подобрено линейно осветление, </w:t>
      </w:r>
      <w:r>
        <w:t>HDR</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t>, работни процеси за скайбокс и кубично картографиране, подобрена система за разпределение на задачите, нов „</w:t>
      </w:r>
      <w:r>
        <w:t>CPU</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t xml:space="preserve"> Таймлайн Профайлър” позволяващ да се следи многоядреното използване, подобрена система за навигиране НавМеш.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t xml:space="preserve">Допреди Юнити 5.0 енджинът използваше доста остаряла версия на физичния мидълуеър </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>PhysX</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t xml:space="preserve"> на </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Nvidia</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t xml:space="preserve">. Юнити 5.0 включва версия 3.3, която е стандарт за </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">Triple-A </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t>игрите.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t xml:space="preserve">Юнити 5.0 осигурява съпорт за Уиндоус, </w:t>
      </w:r>
      <w:r>
        <w:t>OS X</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t xml:space="preserve">, Юнити Уебплейър, Андроид, </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>iOS</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t xml:space="preserve">, Блекбери 10, </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">Windows Phone 8, </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Tizen</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">, </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>WebGL</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">, </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t xml:space="preserve">Плейстейшън 3, Плейстейшън 4, Плейстейшън Вита, </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Wii</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> U, Nintendo 3DS line, Xbox 360, Xbox One, </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t xml:space="preserve">Андроид ТВ, Самсунг Смарт ТВ, </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">Oculus Rift, HTC Vive </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t xml:space="preserve">и </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">     Gear VR.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="bg-BG"/>
        </w:rPr>
        <w:t>На 26-ти Август 2015г. е реализирана експериментална версия за Убунту Линукс с неофициален платформен инсталатор за повечето модерни 64-битови версии на Линукс и официален съпорт само за 64-битово Убунту 12.04 или по-ново.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="PlainText"/>
        <w:rPr>
          <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
        </w:rPr>
      </w:pPr>
    </w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertRange.InsertXML($xml)
Write-Output "done"
Write-Output $d.Paragraphs.Count
